$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B183").Value = 4300
$ws.Range("D183").Value = 137
$ws.Range("B184").Value = 4526
$ws.Range("D184").Value = 226
$ws.Range("B185").Value = 4614
$ws.Range("D185").Value = 88
$ws.Range("B186").Value = 4636
$ws.Range("D186").Value = 22
$ws.Range("B187").Value = 4727
$ws.Range("D187").Value = 91
$ws.Range("B188").Value = 4888
$ws.Range("D188").Value = 161
$ws.Range("B189").Value = 5066
$ws.Range("D189").Value = 178
$ws.Range("B190").Value = 5252
$ws.Range("D190").Value = 186
$ws.Range("B191").Value = 5453
$ws.Range("D191").Value = 201
$ws.Range("B192").Value = 5532
$ws.Range("D192").Value = 79
$ws.Range("B193").Value = 5580
$ws.Range("D193").Value = 48
$ws.Range("B194").Value = 5768
$ws.Range("D194").Value = 188
$ws.Range("B195").Value = 5860
$ws.Range("D195").Value = 92
$ws.Range("B196").Value = 6021
$ws.Range("D196").Value = 161
$ws.Range("B197").Value = 6256
$ws.Range("D197").Value = 235
$ws.Range("B198").Value = 6546
$ws.Range("D198").Value = 290
$ws.Range("B199").Value = 6677
$ws.Range("D199").Value = 131
$ws.Range("B200").Value = 6756
$ws.Range("D200").Value = 79
$ws.Range("B201").Value = 6931
$ws.Range("D201").Value = 175
$ws.Range("B202").Value = 7269
$ws.Range("D202").Value = 338
$ws.Range("B203").Value = 7629
$ws.Range("D203").Value = 360
$ws.Range("B204").Value = 8048
$ws.Range("D204").Value = 419
$ws.Range("B205").Value = 8600
$ws.Range("D205").Value = 552
$ws.Range("B206").Value = 9078
$ws.Range("D206").Value = 478
$ws.Range("B207").Value = 9343
$ws.Range("D207").Value = 265
$ws.Range("B208").Value = 9574
$ws.Range("D208").Value = 231
$ws.Range("B209").Value = 10141
$ws.Range("D209").Value = 567
$ws.Range("B210").Value = 10938
$ws.Range("D210").Value = 797
$ws.Range("B211").Value = 11617
$ws.Range("D211").Value = 679
$ws.Range("B212").Value = 12321
$ws.Range("D212").Value = 704
$ws.Range("B213").Value = 13139
$ws.Range("D213").Value = 818
$ws.Range("B214").Value = 13492
$ws.Range("D214").Value = 353
$ws.Range("B215").Value = 13812
$ws.Range("D215").Value = 320
$ws.Range("B216").Value = 14689
$ws.Range("D216").Value = 877
$ws.Range("B217").Value = 15726
$ws.Range("D217").Value = 1037
$ws.Range("B218").Value = 16910
$ws.Range("D218").Value = 1184
$ws.Range("B219").Value = 18797
$ws.Range("D219").Value = 1887
$ws.Range("B220").Value = 19851
$ws.Range("D220").Value = 1054
$ws.Range("B221").Value = 20355
$ws.Range("D221").Value = 504
$ws.Range("B222").Value = 20886
$ws.Range("D222").Value = 531
$ws.Range("B223").Value = 22296
$ws.Range("D223").Value = 1410
$ws.Range("B224").Value = 24225
$ws.Range("D224").Value = 1929
$ws.Range("B225").Value = 26300
$ws.Range("D225").Value = 2075
$ws.Range("B226").Value = 28268
$ws.Range("D226").Value = 1968
$ws.Range("B227").Value = 29835
$ws.Range("D227").Value = 1567
$ws.Range("B228").Value = 30695
$ws.Range("D228").Value = 860
$ws.Range("B229").Value = 31400
$ws.Range("D229").Value = 705
$ws.Range("B230").Value = 33602
$ws.Range("D230").Value = 2202
$ws.Range("B231").Value = 35330
$ws.Range("D231").Value = 1728
$ws.Range("B232").Value = 37911
$ws.Range("D232").Value = 2581
$ws.Range("B233").Value = 40801
$ws.Range("D233").Value = 2890
$ws.Range("B234").Value = 43843
$ws.Range("D234").Value = 3042
$ws.Range("B235").Value = 45155
$ws.Range("D235").Value = 1312
$ws.Range("B236").Value = 46056
$ws.Range("D236").Value = 901
$ws.Range("B237").Value = 48943
$ws.Range("D237").Value = 2887
$ws.Range("B238").Value = 51728
$ws.Range("D238").Value = 2785
$ws.Range("B239").Value = 55091
$ws.Range("D239").Value = 3363
$ws.Range("B240").Value = 57664
$ws.Range("D240").Value = 2573
$ws.Range("B241").Value = 59946
$ws.Range("D241").Value = 2282
$ws.Range("B242").Value = 61829
$ws.Range("D242").Value = 1883
$ws.Range("B243").Value = 63556
$ws.Range("D243").Value = 1727
$ws.Range("B244").Value = 66772
$ws.Range("D244").Value = 3216
$ws.Range("B245").Value = 68734
$ws.Range("D245").Value = 1962
$ws.Range("B246").Value = 71088
$ws.Range("D246").Value = 2354
$ws.Range("B247").Value = 73667
$ws.Range("D247").Value = 2579
$ws.Range("B248").Value = 75495
$ws.Range("D248").Value = 1828
$ws.Range("B249").Value = 76072
$ws.Range("D249").Value = 577
$ws.Range("B250").Value = 77123
$ws.Range("D250").Value = 1051
$ws.Range("B251").Value = 79181
$ws.Range("D251").Value = 2058
$ws.Range("B252").Value = 81772
$ws.Range("D252").Value = 2591
$ws.Range("B253").Value = 83796
$ws.Range("D253").Value = 2024
$ws.Range("B254").Value = 85567
$ws.Range("D254").Value = 1771
$ws.Range("B255").Value = 86767
$ws.Range("D255").Value = 1200
$ws.Range("B256").Value = 87276
$ws.Range("D256").Value = 509
$ws.Range("B257").Value = 88602
$ws.Range("D257").Value = 1326
$ws.Range("B258").Value = 89913
$ws.Range("D258").Value = 1311
$ws.Range("B259").Value = 91578
$ws.Range("D259").Value = 1665
$ws.Range("B260").Value = 93396
$ws.Range("D260").Value = 1818
$ws.Range("B261").Value = 95257
$ws.Range("D261").Value = 1861
$ws.Range("B262").Value = 96241
$ws.Range("D262").Value = 984
$ws.Range("B263").Value = 96472
$ws.Range("D263").Value = 231
$ws.Range("B264").Value = 97493
$ws.Range("D264").Value = 1021
$ws.Range("B265").Value = 99304
$ws.Range("D265").Value = 1811
$ws.Range("B266").Value = 101257
$ws.Range("D266").Value = 1953
$ws.Range("B267").Value = 103106
$ws.Range("D267").Value = 1849
$ws.Range("B268").Value = 104633
$ws.Range("D268").Value = 1527
$ws.Range("B269").Value = 105733
$ws.Range("D269").Value = 1100
$ws.Range("B270").Value = 105929
$ws.Range("D270").Value = 196
$ws.Range("B271").Value = 107183
$ws.Range("D271").Value = 1254
$ws.Range("B272").Value = 109226
$ws.Range("D272").Value = 2043
$ws.Range("B273").Value = 111208
$ws.Range("D273").Value = 1982
$ws.Range("B274").Value = 113392
$ws.Range("D274").Value = 2184
$ws.Range("B275").Value = 115462
$ws.Range("D275").Value = 2070
$ws.Range("B276").Value = 116731
$ws.Range("D276").Value = 1269
$ws.Range("B277").Value = 117283
$ws.Range("D277").Value = 552
$ws.Range("B278").Value = 119232
$ws.Range("D278").Value = 1949
$ws.Range("B279").Value = 121796
$ws.Range("D279").Value = 2564
$ws.Range("B280").Value = 124921
$ws.Range("D280").Value = 3125
$ws.Range("B281").Value = 127087
$ws.Range("D281").Value = 2166
$ws.Range("B282").Value = 130794
$ws.Range("D282").Value = 3707
$ws.Range("B283").Value = 132984
$ws.Range("D283").Value = 2190
$ws.Range("B284").Value = 133489
$ws.Range("D284").Value = 505
$ws.Range("B285").Value = 135523
$ws.Range("D285").Value = 2034
$ws.Range("B286").Value = 139088
$ws.Range("D286").Value = 3565
$ws.Range("B287").Value = 142133
$ws.Range("D287").Value = 3045
$ws.Range("B288").Value = 146124
$ws.Range("D288").Value = 3991
$ws.Range("B289").Value = 149275
$ws.Range("D289").Value = 3151
$ws.Range("B290").Value = 151336
$ws.Range("D290").Value = 2061
$ws.Range("B291").Value = 152555
$ws.Range("D291").Value = 1219
$ws.Range("B292").Value = 155218
$ws.Range("D292").Value = 2663
$ws.Range("B293").Value = 158905
$ws.Range("D293").Value = 3687
$ws.Range("B294").Value = 161562
$ws.Range("D294").Value = 2657
$ws.Range("B295").Value = 165608
$ws.Range("D295").Value = 4046
$ws.Range("B296").Value = 166649
$ws.Range("D296").Value = 1041
$ws.Range("B297").Value = 167523
$ws.Range("D297").Value = 874
$ws.Range("B298").Value = 168092
$ws.Range("D298").Value = 569
$ws.Range("B299").Value = 170187
$ws.Range("D299").Value = 2095
$ws.Range("B300").Value = 173228
$ws.Range("D300").Value = 3041
$ws.Range("B301").Value = 179543
$ws.Range("D301").Value = 6315
$ws.Range("B302").Value = 184508
$ws.Range("D302").Value = 4965
$ws.Range("B303").Value = 186244
$ws.Range("D303").Value = 1736
$ws.Range("B304").Value = 187463
$ws.Range("D304").Value = 1219
$ws.Range("B305").Value = 188099
$ws.Range("D305").Value = 636
$ws.Range("B306").Value = 191088
$ws.Range("D306").Value = 2989
$ws.Range("B307").Value = 196047
$ws.Range("D307").Value = 4959
$ws.Range("B308").Value = 198184
$ws.Range("D308").Value = 2137
$ws.Range("B309").Value = 201164
$ws.Range("D309").Value = 2980
$ws.Range("B310").Value = 205236
$ws.Range("D310").Value = 4072
$ws.Range("B311").Value = 208209
$ws.Range("D311").Value = 2973
$ws.Range("B312").Value = 209069
$ws.Range("D312").Value = 860
$ws.Range("B313").Value = 211479
$ws.Range("D313").Value = 2410
$ws.Range("B314").Value = 215055
$ws.Range("D314").Value = 3576
$ws.Range("B315").Value = 217978
$ws.Range("D315").Value = 2923
$ws.Range("B316").Value = 220707
$ws.Range("D316").Value = 2729
$ws.Range("B317").Value = 222752
$ws.Range("D317").Value = 2045
$ws.Range("B318").Value = 223325
$ws.Range("D318").Value = 573
$ws.Range("B319").Value = 224385
$ws.Range("D319").Value = 1060
$ws.Range("B320").Value = 226294
$ws.Range("D320").Value = 1909
$ws.Range("B321").Value = 228778
$ws.Range("D321").Value = 2484
$ws.Range("B322").Value = 231242
$ws.Range("D322").Value = 2464
$ws.Range("B323").Value = 233027
$ws.Range("D323").Value = 1785
$ws.Range("B324").Value = 234571
$ws.Range("D324").Value = 1544
$ws.Range("B325").Value = 236476
$ws.Range("D325").Value = 1905
$ws.Range("B326").Value = 237027
$ws.Range("D326").Value = 551
$ws.Range("B327").Value = 238617
$ws.Range("D327").Value = 1590
$ws.Range("B328").Value = 241392
$ws.Range("D328").Value = 2775
$ws.Range("B329").Value = 243427
$ws.Range("D329").Value = 2035
$ws.Range("B330").Value = 246008
$ws.Range("D330").Value = 2581
$ws.Range("B331").Value = 248190
$ws.Range("D331").Value = 2182
$ws.Range("B332").Value = 249913
$ws.Range("D332").Value = 1723
$ws.Range("B333").Value = 250357
$ws.Range("D333").Value = 444
$ws.Range("B334").Value = 252094
$ws.Range("D334").Value = 1737
$ws.Range("B335").Value = 254826
$ws.Range("D335").Value = 2732
$ws.Range("B336").Value = 256903
$ws.Range("D336").Value = 2077
$ws.Range("B337").Value = 259533
$ws.Range("D337").Value = 2630
$ws.Range("B338").Value = 261774
$ws.Range("D338").Value = 2241
$ws.Range("B339").Value = 263326
$ws.Range("D339").Value = 1552
$ws.Range("B340").Value = 264083
$ws.Range("D340").Value = 757
$ws.Range("B341").Value = 265807
$ws.Range("D341").Value = 1724
$ws.Range("B342").Value = 268986
$ws.Range("D342").Value = 3179
$ws.Range("B343").Value = 271473
$ws.Range("D343").Value = 2487
$ws.Range("B344").Value = 273904
$ws.Range("D344").Value = 2431
$ws.Range("B345").Value = 276234
$ws.Range("D345").Value = 2330
$ws.Range("B346").Value = 277682
$ws.Range("D346").Value = 1448
$ws.Range("B347").Value = 278254
$ws.Range("D347").Value = 572
$ws.Range("B348").Value = 279696
$ws.Range("D348").Value = 1442

# Add new row 349
$ws.Range("A349").Value = 44243
$ws.Range("A349").NumberFormat = $ws.Range("A348").NumberFormat
$ws.Range("B349").Value = 282864
$ws.Range("C349").Value = 13050
$ws.Range("D349").Value = 3168
$ws.Range("E349").Value = 6168
$ws.Range("F349").Value = 140661
$ws.Range("G349").Value = 3881
